$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 1050.8572
$ws.Range("J53").Value = 688.375
$ws.Range("L53").Value = 688.375
$ws.Range("N53").Value = -1962.375

# Row 96
$ws.Range("H96").Value = 2290.4546
$ws.Range("I96").Value = 504.625
$ws.Range("K96").Value = 1513.875
$ws.Range("M96").Value = -140.875

# Row 97
$ws.Range("H97").Value = 849.5
$ws.Range("J97").Value = 849.5
$ws.Range("L97").Value = 2548.5
$ws.Range("N97").Value = -3540.5

# Row 116
$ws.Range("H116").Value = 2784233
$ws.Range("I116").Value = 8000
$ws.Range("K116").Value = 8000
$ws.Range("M116").Value = -4558

# Row 123
$ws.Range("H123").Value = 64475.715
$ws.Range("J123").Value = 64475.715
$ws.Range("L123").Value = 64475.715
$ws.Range("N123").Value = -74275.715

# Row 132
$ws.Range("H132").Value = 2369.2632
$ws.Range("I132").Value = 1945.5
$ws.Range("K132").Value = 5836.5
$ws.Range("M132").Value = -3306.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11486.853
$ws.Range("I32").Value = 10654.744
$ws.Range("K32").Value = 10654.744
$ws.Range("M32").Value = -10367.744

# Row 61
$ws.Range("H61").Value = 1278.2
$ws.Range("I61").Value = 610.1429000000001
$ws.Range("J61").Value = 2837
$ws.Range("K61").Value = 610.1429000000001
$ws.Range("L61").Value = 2837
$ws.Range("M61").Value = -398.1429000000001
$ws.Range("N61").Value = -3261

# Row 97
$ws.Range("H97").Value = 684.7143
$ws.Range("I97").Value = 665.8889
$ws.Range("K97").Value = 665.8889
$ws.Range("M97").Value = -169.8889

# Row 136
$ws.Range("H136").Value = 1278.2
$ws.Range("I136").Value = 610.1429000000001
$ws.Range("J136").Value = 2837
$ws.Range("K136").Value = 1830.4287
$ws.Range("L136").Value = 8511
$ws.Range("M136").Value = 719.5712999999998
$ws.Range("N136").Value = -13611

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 235168.8
$ws.Range("I20").Value = 293784.88
$ws.Range("J20").Value = 704.5
$ws.Range("K20").Value = 293784.88
$ws.Range("L20").Value = 704.5
$ws.Range("M20").Value = -293537.88
$ws.Range("N20").Value = -1198.5

# Row 94
$ws.Range("H94").Value = 3131.5417
$ws.Range("I94").Value = 3120.7727
$ws.Range("J94").Value = 3250
$ws.Range("K94").Value = 3120.7727
$ws.Range("L94").Value = 3250
$ws.Range("M94").Value = -2669.7727
$ws.Range("N94").Value = -4152

# Row 132
$ws.Range("H132").Value = 33730.31
$ws.Range("J132").Value = 33730.31
$ws.Range("L132").Value = 33730.31
$ws.Range("N132").Value = -43850.31

# Row 134
$ws.Range("H134").Value = 2986.2092
$ws.Range("I134").Value = 2356.5898
$ws.Range("K134").Value = 7069.769400000001
$ws.Range("M134").Value = -4534.769400000001

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 2818.04
$ws.Range("I122").Value = 2508.1052
$ws.Range("J122").Value = 3799.5
$ws.Range("K122").Value = 7524.3156
$ws.Range("L122").Value = 11398.5
$ws.Range("M122").Value = -5074.3156
$ws.Range("N122").Value = -16298.5

# Row 132
$ws.Range("H132").Value = 3095.25
$ws.Range("I132").Value = 3266.3333
$ws.Range("J132").Value = 2992.6
$ws.Range("K132").Value = 9798.999899999999
$ws.Range("L132").Value = 8977.799999999999
$ws.Range("M132").Value = -7268.999899999999
$ws.Range("N132").Value = -14037.8

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1377.4117
$ws.Range("J131").Value = 1957.5
$ws.Range("L131").Value = 5872.5
$ws.Range("N131").Value = -15952.5

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 21945.834
$ws.Range("J21").Value = 22304.545
$ws.Range("L21").Value = 22304.545
$ws.Range("N21").Value = -22650.545

# Row 30
$ws.Range("H30").Value = 21945.834
$ws.Range("J30").Value = 22304.545
$ws.Range("L30").Value = 22304.545
$ws.Range("N30").Value = -22514.545

# Row 46
$ws.Range("H46").Value = 14910.5
$ws.Range("I46").Value = 3213.6667
$ws.Range("K46").Value = 3213.6667
$ws.Range("M46").Value = -3057.6667

# Row 52
$ws.Range("H52").Value = 68553.14
$ws.Range("J52").Value = 26394.4
$ws.Range("L52").Value = 26394.4
$ws.Range("N52").Value = -26912.4

# Row 57
$ws.Range("H57").Value = 23372.75
$ws.Range("J57").Value = 29330.5
$ws.Range("L57").Value = 29330.5
$ws.Range("N57").Value = -30970.5

# Row 58
$ws.Range("H58").Value = 28748.75
$ws.Range("J58").Value = 28748.75
$ws.Range("L58").Value = 28748.75
$ws.Range("N58").Value = -29302.75

# Row 97
$ws.Range("H97").Value = 950
$ws.Range("I97").Value = 950
$ws.Range("K97").Value = 950
$ws.Range("M97").Value = -454

# Row 122
$ws.Range("H122").Value = 254219.78
$ws.Range("I122").Value = 377349.16
$ws.Range("J122").Value = 7961
$ws.Range("K122").Value = 1132047.48
$ws.Range("L122").Value = 23883
$ws.Range("M122").Value = -1129597.48
$ws.Range("N122").Value = -28783

# Row 123
$ws.Range("H123").Value = 50799.2
$ws.Range("J123").Value = 50799.2
$ws.Range("L123").Value = 50799.2
$ws.Range("N123").Value = -55699.2

# Row 132
$ws.Range("H132").Value = 3196.6943
$ws.Range("J132").Value = 4303.846
$ws.Range("L132").Value = 12911.538
$ws.Range("N132").Value = -17971.538

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 4831.1665
$ws.Range("I16").Value = 4297.4
$ws.Range("K16").Value = 4297.4
$ws.Range("M16").Value = -4127.4

# Row 46
$ws.Range("H46").Value = 2474.7646
$ws.Range("I46").Value = 1620.5714
$ws.Range("J46").Value = 3072.7
$ws.Range("K46").Value = 1620.5714
$ws.Range("L46").Value = 3072.7
$ws.Range("M46").Value = -1432.5714
$ws.Range("N46").Value = -3448.7

# Row 68
$ws.Range("H68").Value = 1053002
$ws.Range("I68").Value = 1053002
$ws.Range("K68").Value = 1053002
$ws.Range("M68").Value = -1052253

# Row 71
$ws.Range("H71").Value = 1053002
$ws.Range("I71").Value = 1053002
$ws.Range("K71").Value = 5265010
$ws.Range("M71").Value = -5261266

# Row 100
$ws.Range("H100").Value = 9088.294
$ws.Range("J100").Value = 9930
$ws.Range("L100").Value = 9930
$ws.Range("N100").Value = -11012

# Row 132
$ws.Range("H132").Value = 1851
$ws.Range("I132").Value = 1572.0454
$ws.Range("K132").Value = 4716.1362
$ws.Range("M132").Value = -2186.1362

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 6582451
$ws.Range("I96").Value = 1339
$ws.Range("J96").Value = 26325788
$ws.Range("K96").Value = 1339
$ws.Range("L96").Value = 26325788
$ws.Range("M96").Value = 34
$ws.Range("N96").Value = -26328534

# Row 107
$ws.Range("H107").Value = 8098
$ws.Range("I107").Value = 12852
$ws.Range("J107").Value = 1985.7142
$ws.Range("K107").Value = 38556
$ws.Range("L107").Value = 5957.142599999999
$ws.Range("M107").Value = -36636
$ws.Range("N107").Value = -9797.142599999999

# Row 132
$ws.Range("H132").Value = 1159.9231
$ws.Range("I132").Value = 941.65216
$ws.Range("J132").Value = 2833.3333
$ws.Range("K132").Value = 2824.95648
$ws.Range("L132").Value = 8499.999899999999
$ws.Range("M132").Value = -294.9564799999998
$ws.Range("N132").Value = -13559.9999
